# Updates cryptos list (price + 1h volume change columns) per the
# Wed May 31 11:11:26 UTC 2023 GitHub Actions data refresh.
#
# Column D ("Price") holds values typed as plain text in the source sheet
# (e.g. "27.103.94", "1.002"). Most of the new prices still look numeric to
# Excel ("1.001", "306.99", ...), so a bare `.Value =` assignment would be
# silently reinterpreted as a Number, truncating things like trailing
# zeroes. For those cells we momentarily switch the cell to Text format
# ("@") before writing the literal string, then restore the cell's style
# to "Normal" so no stray number-format/style id is left behind on save.
# Column E ("Volume(1h)") values always carry padding spaces + a trailing
# "%" so Excel already stores them as text with a plain `.Value =`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a Number.
$textCells = @(
    @{ Cell = 'D5'; Value = '306.99' }
    @{ Cell = 'D7'; Value = '0.5067' }
    @{ Cell = 'D8'; Value = '0.3739' }
    @{ Cell = 'D9'; Value = '0.07150' }
    @{ Cell = 'D10'; Value = '0.8860' }
    @{ Cell = 'D11'; Value = '20.58' }
    @{ Cell = 'D13'; Value = '0.07550' }
    @{ Cell = 'D14'; Value = '5.312' }
    @{ Cell = 'D16'; Value = '1.001' }
    @{ Cell = 'D17'; Value = '0.000008467' }
    @{ Cell = 'D18'; Value = '14.10' }
    @{ Cell = 'D21'; Value = '5.057' }
    @{ Cell = 'D23'; Value = '10.56' }
    @{ Cell = 'D24'; Value = '6.485' }
    @{ Cell = 'D25'; Value = '150.77' }
    @{ Cell = 'D26'; Value = '1.839' }
    @{ Cell = 'D27'; Value = '17.92' }
    @{ Cell = 'D28'; Value = '2.099' }
    @{ Cell = 'D29'; Value = '112.63' }
    @{ Cell = 'D30'; Value = '4.740' }
    @{ Cell = 'D31'; Value = '4.677' }
    @{ Cell = 'D32'; Value = '0.09035' }
    @{ Cell = 'D34'; Value = '3.098' }
    @{ Cell = 'D35'; Value = '1.158' }
    @{ Cell = 'D36'; Value = '0.7368' }
    @{ Cell = 'D37'; Value = '0.02036' }
    @{ Cell = 'D38'; Value = '2.492' }
    @{ Cell = 'D39'; Value = '3.043' }
    @{ Cell = 'D40'; Value = '1.081' }
    @{ Cell = 'D41'; Value = '0.5335' }
    @{ Cell = 'D43'; Value = '115.39' }
    @{ Cell = 'D44'; Value = '8.336' }
    @{ Cell = 'D45'; Value = '0.1471' }
    @{ Cell = 'D46'; Value = '0.4637' }
    @{ Cell = 'D47'; Value = '1.000' }
    @{ Cell = 'D48'; Value = '10.01' }
    @{ Cell = 'D49'; Value = '1.563' }
    @{ Cell = 'D50'; Value = '64.55' }
    @{ Cell = 'D51'; Value = '36.32' }
)

# Cells that are safe to assign directly (non-numeric-looking strings:
# coin names, links, the multi-dot big-BTC-style prices, and the
# percentage-change column).
$plainCells = @(
    @{ Cell = 'D2'; Value = '27.038.58' }
    @{ Cell = 'E2'; Value = '  -3.15%  ' }
    @{ Cell = 'D3'; Value = '1.865.54' }
    @{ Cell = 'E3'; Value = '  -2.27%  ' }
    @{ Cell = 'E4'; Value = '  +0.30%  ' }
    @{ Cell = 'E5'; Value = '  -2.06%  ' }
    @{ Cell = 'E6'; Value = '  +0.23%  ' }
    @{ Cell = 'E7'; Value = '  +0.96%  ' }
    @{ Cell = 'E9'; Value = '  -2.38%  ' }
    @{ Cell = 'E10'; Value = '  -2.81%  ' }
    @{ Cell = 'E11'; Value = '  -3.09%  ' }
    @{ Cell = 'D12'; Value = '1.866.89' }
    @{ Cell = 'E12'; Value = '  -2.89%  ' }
    @{ Cell = 'E13'; Value = '  -1.59%  ' }
    @{ Cell = 'E14'; Value = '  -3.21%  ' }
    @{ Cell = 'E15'; Value = '  -4.00%  ' }
    @{ Cell = 'E16'; Value = '  +0.16%  ' }
    @{ Cell = 'E17'; Value = '  -3.21%  ' }
    @{ Cell = 'E19'; Value = '  +0.05%  ' }
    @{ Cell = 'D20'; Value = '27.121.98' }
    @{ Cell = 'E20'; Value = '  -2.97%  ' }
    @{ Cell = 'E21'; Value = '  -2.46%  ' }
    @{ Cell = 'D22'; Value = '2.111.96' }
    @{ Cell = 'E22'; Value = '  -1.98%  ' }
    @{ Cell = 'E23'; Value = '  -2.67%  ' }
    @{ Cell = 'E24'; Value = '  -1.99%  ' }
    @{ Cell = 'E25'; Value = '  -1.51%  ' }
    @{ Cell = 'E26'; Value = '  -0.40%  ' }
    @{ Cell = 'E27'; Value = '  -2.74%  ' }
    @{ Cell = 'E28'; Value = '  -4.83%  ' }
    @{ Cell = 'E29'; Value = '  -2.56%  ' }
    @{ Cell = 'E30'; Value = '  -3.79%  ' }
    @{ Cell = 'E31'; Value = '  -3.97%  ' }
    @{ Cell = 'E32'; Value = '  -0.05%  ' }
    @{ Cell = 'E33'; Value = '  -3.23%  ' }
    @{ Cell = 'E34'; Value = '  -3.54%  ' }
    @{ Cell = 'E35'; Value = '  -6.46%  ' }
    @{ Cell = 'E36'; Value = '  -5.27%  ' }
    @{ Cell = 'E37'; Value = '  -2.30%  ' }
    @{ Cell = 'E38'; Value = '  -3.72%  ' }
    @{ Cell = 'E39'; Value = '  -0.59%  ' }
    @{ Cell = 'E40'; Value = '  -1.31%  ' }
    @{ Cell = 'E41'; Value = '  -4.00%  ' }
    @{ Cell = 'E42'; Value = '  -4.25%  ' }
    @{ Cell = 'E43'; Value = '  +1.68%  ' }
    @{ Cell = 'E44'; Value = '  -2.20%  ' }
    @{ Cell = 'E45'; Value = '  -3.42%  ' }
    @{ Cell = 'B46'; Value = 'Decentraland' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'E46'; Value = '  -4.14%  ' }
    @{ Cell = 'B47'; Value = 'PaxDollar' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Cell = 'E47'; Value = '  +0.19%  ' }
    @{ Cell = 'E48'; Value = '  -5.80%  ' }
    @{ Cell = 'E49'; Value = '  -4.60%  ' }
    @{ Cell = 'E50'; Value = '  -4.55%  ' }
    @{ Cell = 'E51'; Value = '  -2.34%  ' }
)

foreach ($item in $textCells) {
    $cell = $ws.Range($item.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

foreach ($item in $plainCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
